$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Notified Production (MW)" values for rows 2..97 (quarterly data for the new day,
# fixing the huge difference between the hourly and quarterly forecasted values)
$values = @(26,26,26,26,25,25,25,25,24,24,24,23,23,23,23,23,35,34,35,35,32,31,30,30,25,25,25,25,56,56,56,57,66,66,67,69,56,57,58,61,87,88,90,90,85,86,86,87,110,110,111,112,140,141,142,142,203,203,203,204,237,237,238,239,304,303,304,304,379,380,381,382,406,408,410,411,391,391,392,392,368,368,367,366,333,332,331,332,285,284,283,282,0,0,0,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    # Shift the timestamp forward by 3 days, preserving the time-of-day fraction
    $oldDate = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 1).Value = $oldDate + 3
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
